$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: find a shape by (partial) name within a Shapes collection.
# ---------------------------------------------------------------------------
function Find-ShapeByName($shapes, $namePrefix) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "$namePrefix*") {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" field text (3/23/2018 -> 3/24/2018)
#    on every slide layout, the slide master, and the notes master.
# ---------------------------------------------------------------------------
$oldDate = "3/23/2018"
$newDate = "3/24/2018"

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($j)
    $dateShape = Find-ShapeByName $layout.Shapes "Date Placeholder"
    if ($dateShape -ne $null) {
        $dateShape.TextFrame.TextRange.Text = $newDate
    }
}

$masterDateShape = Find-ShapeByName $p.SlideMaster.Shapes "Date Placeholder"
if ($masterDateShape -ne $null) {
    $masterDateShape.TextFrame.TextRange.Text = $newDate
}

$notesDateShape = Find-ShapeByName $p.NotesMaster.Shapes "Date Placeholder"
if ($notesDateShape -ne $null) {
    $notesDateShape.TextFrame.TextRange.Text = $newDate
}

# ---------------------------------------------------------------------------
# 2) Slide 21 ("DEMO"): retitle the first bullet and repoint the link.
# ---------------------------------------------------------------------------
$slide21 = $p.Slides.Item(21)
$demoShape = $slide21.Shapes.Item(2)
$demoRange = $demoShape.TextFrame.TextRange

$oldTitle = "Fibonacci sequence"
$newTitle = "Even-odd"
$oldUrl = "https://github.com/shankar-ray/Assembly-Language-Tutorials-for-Windows/blob/master/03%20Assembly%20Language/Assembly%20Language/Assembly%20Language/Source.asm"
$newUrl = "https://github.com/shankar-ray/Assembly-Language-Tutorials-for-Windows/tree/master/03%20Assembly%20Language"

# Replace the URL first so the earlier offsets (for the title) stay valid.
$urlStart = $demoRange.Length - $oldUrl.Length + 1
$demoRange.Characters($urlStart, $oldUrl.Length).Text = $newUrl
$demoRange.Characters(1, $oldTitle.Length).Text = $newTitle

# ---------------------------------------------------------------------------
# 3) Slide 8 ("IDENTIFIERS" / "DIRECTIVES"): fix the un-initialized data typo
#    and collapse the three runs describing it into a single run.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$directivesShape = $slide8.Shapes.Item(4)
$directivesRange = $directivesShape.TextFrame.TextRange

$oldParagraph = ".data? //segment contains un-initalized data"
$newParagraph = ".data? //segment contains un-initialized data"

$fullText = $directivesRange.Text
$paraStart = $fullText.IndexOf($oldParagraph) + 1
$directivesRange.Characters($paraStart, $oldParagraph.Length).Text = $newParagraph
